$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-shape the row layout.
#    The lone "    " spacer row (old row 12) is replaced by nine new rows
#    (three new timesheet entries + their git-commit separators, plus three
#    blank spacer rows) which pushes everything from the old row 13 onward
#    down by 8 rows (old13 -> new21, old18 -> new26, old26 -> new34).
# ---------------------------------------------------------------------------
$ws.Rows("12:12").Delete()
$ws.Rows("12:20").Insert()

# ---------------------------------------------------------------------------
# 2. New timesheet entry rows (2/8/2014 work session split into 3 chunks)
# ---------------------------------------------------------------------------

# Row 12 - 8:20a-9:20a, Added tests for mouse events
$ws.Range("A12").NumberFormat = "m/d;@"
$ws.Range("A12").Value = 41678
$ws.Range("B12").Value = "Doyle"
$ws.Range("C12").NumberFormat = "h:mm"
$ws.Range("C12").Value = 500/1440
$ws.Range("D12").NumberFormat = "h:mm"
$ws.Range("D12").Value = 560/1440
$ws.Range("E12").Value = 0
$ws.Range("G12").Value = 2
$ws.Range("I12").Formula = "=60"
$ws.Range("K12").Value = "Added tests for mouse events"

# Row 13 - commit separator
$ws.Range("A13").Value = " =========================committing to git: 2/8 9:21 ========================="

# Row 14 - 9:45a-11:12a, More tests for mouse and key events
$ws.Range("A14").NumberFormat = "m/d;@"
$ws.Range("A14").Value = 41678
$ws.Range("B14").Value = "Doyle"
$ws.Range("C14").NumberFormat = "h:mm"
$ws.Range("C14").Value = 585/1440
$ws.Range("D14").NumberFormat = "h:mm"
$ws.Range("D14").Value = 672/1440
$ws.Range("E14").Value = 1
$ws.Range("G14").Value = 2
$ws.Range("I14").Formula = "=60+12"
$ws.Range("K14").Value = "More tests for mouse and key events"

# Row 15 - commit separator
$ws.Range("A15").Value = " =========================committing to git: 2/8 11:14 ========================="

# Row 16 - 11:43a-1:03p, Got drag working better
$ws.Range("A16").NumberFormat = "m/d;@"
$ws.Range("A16").Value = 41678
$ws.Range("B16").Value = "Doyle"
$ws.Range("C16").NumberFormat = "h:mm"
$ws.Range("C16").Value = 703/1440
$ws.Range("D16").NumberFormat = "h:mm"
$ws.Range("D16").Value = 783/1440
$ws.Range("E16").Value = 1
$ws.Range("G16").Value = 2
$ws.Range("I16").Formula = "=60+10"
$ws.Range("K16").Value = "Got drag working better"

# Row 17 - commit separator
$ws.Range("A17").Value = " =========================committing to git: 2/8 13:15 ========================="

# Rows 18-20 stay blank (spacer rows before the totals section); row 18 keeps
# formatted-but-empty A/C/D cells matching the rest of the time columns.
$ws.Range("C18").NumberFormat = "h:mm"
$ws.Range("D18").NumberFormat = "h:mm"

# ---------------------------------------------------------------------------
# 3. Fix up the totals-section formulas that moved down with the insert but
#    need their ranges/values corrected to match the new data extent.
# ---------------------------------------------------------------------------
$ws.Range("G21").Formula = "=SUMIF(G2:G20,""1"",I2:I20)"
$ws.Range("G22").Formula = "=SUMIF(G3:G20,""2"",I3:I20)"
$ws.Range("G23").Formula = "=SUMIF(G3:G20,""3"",I3:I20)"
$ws.Range("G24").Formula = "=SUMIF(G3:G20,""4"",I3:I20)"

$ws.Range("G26").Formula = "=G13/60"
$ws.Range("G27").Formula = "=G22/60"
$ws.Range("G28").Formula = "=G23/60"
$ws.Range("G29").Formula = "=G24/60"

Write-Output "done"
